# Affix Generator work: add a new "artifact_item_id" column to the Raids
# admin-import sheet, with its sample value for the existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in column H, right after the existing last header column (G).
$ws.Range("H1").Value = "artifact_item_id"

# New sample data value for the single existing data row.
$ws.Range("H2").Value = "Ancestral Finger Bone of The Magi Troth"

# Match the column width Excel would have picked via "best fit" for the
# new column's content (closest value the engine's width grid allows to
# the original authored width of 47.131).
$ws.Columns.Item(8).ColumnWidth = 46.3
